$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for dates/quality/volume/prices/units/price-per-kg got
# shuffled between rows while the market/product/category/origin columns
# stayed put. Apply the new values per row as described by the diff.

# Row 2 (was D2=44475/Especial/...) -> becomes the former Row 6 data
$ws.Range("D2").Value = 44441
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 29000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 29500
$ws.Range("S2").Value = 2458
$ws.Range("T2").Value = 12

# Row 3 -> becomes the former Row 7 data
$ws.Range("D3").Value = 44496
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 23500
$ws.Range("Q3").Value = "$/caja 12 kilos"
$ws.Range("S3").Value = 1958
$ws.Range("T3").Value = 12

# Row 4 -> becomes the former Row 5 data
$ws.Range("D4").Value = 44524
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 23000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 23500
$ws.Range("S4").Value = 1958
$ws.Range("T4").Value = 12

# Row 5 -> becomes the former Row 4 data
$ws.Range("D5").Value = 44482
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 26000
$ws.Range("P5").Value = 25500
$ws.Range("S5").Value = 2125
$ws.Range("T5").Value = 12

# Row 6 -> becomes the former Row 2 data
$ws.Range("D6").Value = 44475
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 32000
$ws.Range("O6").Value = 33000
$ws.Range("P6").Value = 32500
$ws.Range("S6").Value = 2708
$ws.Range("T6").Value = 12

# Row 7 -> becomes the former Row 3 data
$ws.Range("D7").Value = 44468
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("S7").Value = 2950
$ws.Range("T7").Value = 10

# Row 9 -> becomes the former Row 10 data
$ws.Range("D9").Value = 44167
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("S9").Value = 1423

# Row 10 -> becomes the former Row 9 data
$ws.Range("D10").Value = 44160
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("S10").Value = 1500
